# Frappy_Project_Plan.docx edit:
#   1. Change the "Severity - High" risk rating of the "Database Structure"
#      risk item to "Severity - Medium".
#   2. Add a new "Tools - ..." paragraph describing the project's tech
#      stack, right after the existing "Github - ..." paragraph and before
#      the "Configuration Management" heading.

$d = $word.ActiveDocument
$paras = $d.Paragraphs

# ---------------------------------------------------------------------
# Step 1: "Severity - High" -> "Severity - Medium" for Database Structure.
# There are two "Severity - High" paragraphs in the doc (User
# Authentication, then Database Structure) - we want the second one, which
# immediately follows the "Database Structure" / "Likelihood - Low" pair.
# ---------------------------------------------------------------------
$severityHighSeen = 0
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($txt -eq "Severity - High") {
        $severityHighSeen += 1
        if ($severityHighSeen -eq 2) {
            $targetIndex = $i
        }
    }
}

if ($targetIndex -ne -1) {
    $paras.Item($targetIndex).Range.Text = "Severity - Medium"
} else {
    Write-Output "WARNING: could not locate second 'Severity - High' paragraph"
}

# ---------------------------------------------------------------------
# Step 2: Insert the new "Tools - ..." paragraph right after "Github - ...".
# ---------------------------------------------------------------------
$githubIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($txt -like "Github - Website for hosting the git repository*") {
        $githubIndex = $i
        break
    }
}

if ($githubIndex -ne -1) {
    $githubPara = $d.Paragraphs.Item($githubIndex)
    $githubPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($githubIndex + 1)
    $newPara.Range.Text = "Tools - The project will be using a custom stack consisting of React, Django, and Postgres (RPD), where React provides a front end framework for building web and mobile applications and allows us to fetch page specific javascript allowing for a simple to use SPA. Django provides the REST API, user authentication, and serves static files via a reverse proxy with NGINX and postgres is a robust relational database"
} else {
    Write-Output "WARNING: could not locate 'Github - ...' paragraph"
}

Write-Output "Done: severityIndex=$targetIndex githubIndex=$githubIndex"
